$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '40.275.94'
$ws.Range("E2").Value = '  +0.11%  '

# Row 3
$ws.Range("D3").Value = '2.237.79'
$ws.Range("E3").Value = '  -0.49%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").Value = '''294.52'
$ws.Range("E5").Value = '  -0.10%  '

# Row 6
$ws.Range("D6").Value = '''88.86'
$ws.Range("E6").Value = '  +1.88%  '

# Row 7
$ws.Range("E7").Value = '  -0.90%  '

# Row 8
$ws.Range("E8").Value = '  -0.01%  '

# Row 9
$ws.Range("D9").Value = '''0.475'
$ws.Range("E9").Value = '  -0.38%  '

# Row 10
$ws.Range("D10").Value = '''30.43'
$ws.Range("E10").Value = '  -3.12%  '

# Row 11
$ws.Range("D11").Value = '''0.0784'
$ws.Range("E11").Value = '  -2.03%  '

# Row 12
$ws.Range("D12").Value = '''0.112'
$ws.Range("E12").Value = '  +3.07%  '

# Row 13
$ws.Range("D13").Value = '''6.53'
$ws.Range("E13").Value = '  +0.35%  '

# Row 14
$ws.Range("D14").Value = '2.583.30'
$ws.Range("E14").Value = '  -0.37%  '

# Row 15
$ws.Range("D15").Value = '''13.88'
$ws.Range("E15").Value = '  -3.03%  '

# Row 16
$ws.Range("D16").Value = '2.219.79'
$ws.Range("E16").Value = '  -0.43%  '

# Row 17
$ws.Range("D17").Value = '''0.736'
$ws.Range("E17").Value = '  -0.79%  '

# Row 18
$ws.Range("D18").Value = '40.192.85'

# Row 19
$ws.Range("E19").Value = '  -0.58%  '

# Row 20
$ws.Range("E20").Value = '  +6.76%  '

# Row 21
$ws.Range("E21").Value = '  -0.27%  '

# Row 22
$ws.Range("D22").Value = '''65.68'
$ws.Range("E22").Value = '  -0.23%  '

# Row 23
$ws.Range("D23").Value = '''237.25'
$ws.Range("E23").Value = '  +0.18%  '

# Row 24
$ws.Range("E24").Value = '  -0.03%  '

# Row 25
$ws.Range("E25").Value = '  -0.40%  '

# Row 26
$ws.Range("E26").Value = '  -1.79%  '

# Row 27
$ws.Range("D27").Value = '''22.83'
$ws.Range("E27").Value = '  -1.18%  '

# Row 28
$ws.Range("D28").Value = '''2.21'
$ws.Range("E28").Value = '  -0.82%  '

# Row 29
$ws.Range("D29").Value = '''9.27'
$ws.Range("E29").Value = '  -0.48%  '

# Row 30
$ws.Range("D30").Value = '''155.59'
$ws.Range("E30").Value = '  +1.42%  '

# Row 31
$ws.Range("D31").Value = '''32.29'
$ws.Range("E31").Value = '  -3.66%  '

# Row 32
$ws.Range("E32").Value = '  -0.06%  '

# Row 33
$ws.Range("D33").Value = '''4.96'
$ws.Range("E33").Value = '  +0.50%  '

# Row 34
$ws.Range("D34").Value = '''0.0719'
$ws.Range("E34").Value = '  -0.19%  '

# Row 35
$ws.Range("E35").Value = '  -1.47%  '

# Row 36
$ws.Range("E36").Value = '  +5.98%  '

# Row 37
$ws.Range("E37").Value = '  +0.18%  '

# Row 38
$ws.Range("D38").Value = '''15.87'
$ws.Range("E38").Value = '  -5.35%  '

# Row 39
$ws.Range("D39").Value = '''0.0979'
$ws.Range("E39").Value = '  -3.57%  '

# Row 40
$ws.Range("E40").Value = '  -0.69%  '

# Row 41
$ws.Range("D41").Value = '2.140.35'
$ws.Range("E41").Value = '  +5.59%  '

# Row 42
$ws.Range("D42").Value = '''3.87'
$ws.Range("E42").Value = '  +0.80%  '

# Row 43
$ws.Range("D43").Value = '''18.23'
$ws.Range("E43").Value = '  +10.82%  '

# Row 44
$ws.Range("D44").Value = '''2.14'
$ws.Range("E44").Value = '  -4.35%  '

# Row 45
$ws.Range("D45").Value = '''0.0269'

# Row 46
$ws.Range("D46").Value = '''9.85'
$ws.Range("E46").Value = '  -1.15%  '

# Row 47
$ws.Range("D47").Value = '''2.71'
$ws.Range("E47").Value = '  +4.34%  '

# Row 48
$ws.Range("D48").Value = '2.448.71'
$ws.Range("E48").Value = '  -1.04%  '

# Row 49
$ws.Range("E49").Value = '  +1.42%  '

# Row 50
$ws.Range("B50").Value = 'BitcoinSV'
$ws.Range("C50").Value = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range("D50").Value = '''69.64'
$ws.Range("E50").Value = '  -3.23%  '

# Row 51
$ws.Range("B51").Value = 'TrustWalletToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D51").Value = '''1.10'
$ws.Range("E51").Value = '  +0.48%  '
